$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 321, pushing the existing row 321..400 down to 322..401
$ws.Rows(321).Insert()

# Populate the newly inserted row 321 with the new record
$ws.Cells.Item(321, 1).Value = 4
$ws.Cells.Item(321, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(321, 3).Value = "Los Lagos"
$ws.Cells.Item(321, 4).Value = 44964
$ws.Cells.Item(321, 5).Value = 10
$ws.Cells.Item(321, 6).Value = 100112040
$ws.Cells.Item(321, 7).Value = "Cilantro"
$ws.Cells.Item(321, 8).Value = "Sin especificar"
$ws.Cells.Item(321, 9).Value = "Primera"
$ws.Cells.Item(321, 10).Value = 120
$ws.Cells.Item(321, 11).Value = 8000
$ws.Cells.Item(321, 12).Value = 8000
$ws.Cells.Item(321, 13).Value = 8000
$ws.Cells.Item(321, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(321, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(321, 16).Value = 4000
$ws.Cells.Item(321, 17).Value = 2
$ws.Cells.Item(321, 18).Value = "Hortaliza"
